# 678-MS-EPP-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment-Makerepayment1.xlsx
# The early-repayment transaction date moved earlier (see "Input" sheet,
# repaymenttransactiondate = 09-Jan-2015 / serial 42019), so the whole
# amortization schedule (Repayment Schedule), its Summary roll-up and the
# Transactions log were all recalculated/re-entered with the new figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B2").Value = 838.56
$summary.Range("E2").Value = 9161.44
$summary.Range("F2").Value = 930.02

$summary.Range("A3").Value = 641.35
$summary.Range("B3").Value = 96.69
$summary.Range("E3").Value = 544.66
$summary.Range("F3").Value = 84.34

$summary.Range("C4").Select()

# ---------------------------------------------------------------------
# Repayment Schedule sheet
# ---------------------------------------------------------------------
$sched = $wb.Worksheets.Item("Repayment Schedule")

# D3 ("Paid Date") is no longer filled in (the repayment row now reflects
# the brand-new, earlier transaction) - clear it and keep its format in
# step with the adjoining (already blank) E3 cell.
$sched.Range("D3").ClearContents()
$sched.Range("F3").Copy()
$sched.Range("D3:E3").PasteSpecial(-4122)

$sched.Range("F3").Value = 935.25
$sched.Range("G3").Value = 9064.75
$sched.Range("H3").Value = 96.69

# K3 ("Due") crosses the 1,000 mark, so pick up the thousands-formatted
# look already used by the "Balance of Loan" column.
$sched.Range("G3").Copy()
$sched.Range("K3").PasteSpecial(-4122)
$sched.Range("K3").Value = 1031.94

$sched.Range("M3").Value = 935.25
$sched.Range("P3").Value = 96.69

$sched.Range("G4").Value = 8231.42
$sched.Range("H4").Value = 84.34
$sched.Range("K4").Value = 917.67
$sched.Range("P4").Value = 917.67

$sched.Range("G5").Value = 7398.09
$sched.Range("H5").Value = 93.37
$sched.Range("K5").Value = 926.7
$sched.Range("P5").Value = 926.7

$sched.Range("G6").Value = 6564.76
$sched.Range("H6").Value = 72.97
$sched.Range("K6").Value = 906.3
$sched.Range("P6").Value = 906.3

$sched.Range("G7").Value = 5731.43
$sched.Range("H7").Value = 66.91
$sched.Range("K7").Value = 900.24
$sched.Range("P7").Value = 900.24

$sched.Range("G8").Value = 4898.1000000000004
$sched.Range("H8").Value = 56.53
$sched.Range("K8").Value = 889.86
$sched.Range("P8").Value = 889.86

$sched.Range("G9").Value = 4064.77
$sched.Range("H9").Value = 49.92
$sched.Range("K9").Value = 883.25
$sched.Range("P9").Value = 883.25

$sched.Range("G10").Value = 3231.44
$sched.Range("H10").Value = 41.43
$sched.Range("K10").Value = 874.76
$sched.Range("P10").Value = 874.76

$sched.Range("G11").Value = 2398.11
$sched.Range("H11").Value = 31.87
$sched.Range("K11").Value = 865.2
$sched.Range("P11").Value = 865.2

$sched.Range("G12").Value = 1564.78
$sched.Range("H12").Value = 24.44
$sched.Range("K12").Value = 857.77
$sched.Range("P12").Value = 857.77

$sched.Range("G13").Value = 731.45
$sched.Range("H13").Value = 15.43
$sched.Range("K13").Value = 848.76
$sched.Range("P13").Value = 848.76

$sched.Range("F14").Value = 731.45
$sched.Range("H14").Value = 7.45
$sched.Range("K14").Value = 738.9
$sched.Range("P14").Value = 738.9

$sched.Range("C11").Select()

# ---------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------
$txn = $wb.Worksheets.Item("Transactions")

$txn.Range("A2").Value = 35
$txn.Range("C2").Value = 42019
$txn.Range("F2").Value = 838.56
$txn.Range("G2").Value = 96.69
$txn.Range("J2").Value = 9161.44

$txn.Range("A3").Value = 33

# Keep the Transactions tab the active one, matching the workbook's
# original saved state (it was the visible tab before the edit too).
$txn.Range("B2").Select()
